# Weekly price-sheet refresh: a new week's record is inserted at the top of
# the data (row 3, right after the most-recent-anchor row 2), every
# subsequent record shifts down one row, and the oldest record (old row 27)
# is appended as the new last row (28).
#
# Concretely: for each row r in 4..28, the "Fecha" (D) and the measurement
# columns (J..R) take on the values that used to live in row (r-1). Row 3's
# "Fecha" becomes the new week's date (44496) while its own J..R values are
# left untouched (they already hold the correct data, since they are not
# shifted). The sheet grows from 27 to 28 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 3
$lastDataRow = 27
$newLastDataRow = 28
$newDate = 44496

$colD = 4
$colFirstShift = 10   # J
$colLastShift = 18    # R

# 1) Snapshot the current D and J:R values for every data row (3..27)
#    before anything is overwritten.
$oldD = @{}
$oldShift = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, $colD).Value2
    $rowVals = @{}
    for ($c = $colFirstShift; $c -le $colLastShift; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $oldShift[$r] = $rowVals
}

# 2) Write the new last row (28) with the data that used to be row 27's
#    "fixed" columns (A:I) plus row 27's D / J:R values.
for ($c = 1; $c -le 9; $c++) {
    $ws.Cells.Item($newLastDataRow, $c).Value = $ws.Cells.Item($lastDataRow, $c).Value2
}
$ws.Cells.Item($newLastDataRow, $colD).Value = $oldD[$lastDataRow]
$ws.Cells.Item($newLastDataRow, $colD).NumberFormat = $ws.Cells.Item($lastDataRow, $colD).NumberFormat
foreach ($c in $oldShift[$lastDataRow].Keys) {
    $ws.Cells.Item($newLastDataRow, $c).Value = $oldShift[$lastDataRow][$c]
}

# 3) Shift rows 27 down to 4: row r gets row (r-1)'s D and J:R values.
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, $colD).Value = $oldD[$src]
    foreach ($c in $oldShift[$src].Keys) {
        $ws.Cells.Item($r, $c).Value = $oldShift[$src][$c]
    }
}

# 4) Row 3 gets the new week's date; its J:R values are unchanged (still the
#    originally-read row-3 data, never overwritten above).
$ws.Cells.Item($firstDataRow, $colD).Value = $newDate
